$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header strings from _old/_new suffixes to _FV2310/_FV2404 suffixes.
$oldHeaders = @("Segmentname_old","Segmentgruppe_old","Segment_old","Datenelement_old","Segment ID_old","Code_old","Qualifier_old","Beschreibung_old","Bedingungsausdruck_old","Bedingung_old")
$newHeaders = @("Segmentname_new","Segmentgruppe_new","Segment_new","Datenelement_new","Segment ID_new","Code_new","Qualifier_new","Beschreibung_new","Bedingungsausdruck_new","Bedingung_new")

for ($i = 0; $i -lt $oldHeaders.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = ($oldHeaders[$i] -replace "_old$", "_FV2310")
}

for ($i = 0; $i -lt $newHeaders.Length; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = ($newHeaders[$i] -replace "_new$", "_FV2404")
}

# Add a table (ListObject) spanning the used range, with autofilter.
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U91"), [Type]::Missing, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# Freeze the header row.
$ws.Activate()
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

$wb.Save()
